# Insert two new data rows (29 and 30) into the "Ciruela" sheet, pushing
# the existing rows 29..102 down to 31..104.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 2 new rows starting at row 29 (shifts rows 29:102 down to 31:104)
$ws.Rows("29:30").Insert()

# Fill in the two newly-inserted rows with their data.
# Row 29
$ws.Cells.Item(29, 1).Value = 4
$ws.Cells.Item(29, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(29, 3).Value = "Los Lagos"
$ws.Range("D29").Value = 44560
$ws.Cells.Item(29, 5).Value = 10
$ws.Cells.Item(29, 6).Value = "Fruta"
$ws.Cells.Item(29, 7).Value = 100103
$ws.Cells.Item(29, 8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item(29, 9).Value = 100103002
$ws.Cells.Item(29, 10).Value = "Ciruela"
$ws.Cells.Item(29, 11).Value = "Red Beaut"
$ws.Cells.Item(29, 12).Value = "Primera"
$ws.Cells.Item(29, 13).Value = 400
$ws.Cells.Item(29, 14).Value = 18000
$ws.Cells.Item(29, 15).Value = 19000
$ws.Cells.Item(29, 16).Value = 18500
$ws.Cells.Item(29, 17).Value = "`$/caja 15 kilos granel"
$ws.Cells.Item(29, 18).Value = "Región Metropolitana"
$ws.Cells.Item(29, 19).Value = 1233
$ws.Cells.Item(29, 20).Value = 15

# Row 30
$ws.Cells.Item(30, 1).Value = 4
$ws.Cells.Item(30, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(30, 3).Value = "Los Lagos"
$ws.Range("D30").Value = 44560
$ws.Cells.Item(30, 5).Value = 10
$ws.Cells.Item(30, 6).Value = "Fruta"
$ws.Cells.Item(30, 7).Value = 100103
$ws.Cells.Item(30, 8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item(30, 9).Value = 100103002
$ws.Cells.Item(30, 10).Value = "Ciruela"
$ws.Cells.Item(30, 11).Value = "Red Beaut"
$ws.Cells.Item(30, 12).Value = "Segunda"
$ws.Cells.Item(30, 13).Value = 200
$ws.Cells.Item(30, 14).Value = 16000
$ws.Cells.Item(30, 15).Value = 16000
$ws.Cells.Item(30, 16).Value = 16000
$ws.Cells.Item(30, 17).Value = "`$/caja 15 kilos granel"
$ws.Cells.Item(30, 18).Value = "Región Metropolitana"
$ws.Cells.Item(30, 19).Value = 1067
$ws.Cells.Item(30, 20).Value = 15
